# Update "tab_aggiornamenti_alla_documentazione_tecnica.xlsx":
# rename the sheet to reflect the new update date and append the rows
# describing the newly published ANPR documentation updates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet tab to the new update date.
$ws.Name = "aggiornamenti_14_11_2017"

# 2) Add the new log rows (97-105). Each new row re-uses the formatting of
#    an existing, equivalently-styled row so no spurious new cell styles get
#    created; then the real values are written on top of the copied format.
#    NOTE: call these positionally -- named-parameter invocation on this
#    runtime is catastrophically slow when the function body makes COM
#    calls, so every call below passes arguments by position.
function Copy-RowFormat($TemplateRow, $TargetRow) {
    $ws.Range("A$TemplateRow`:D$TemplateRow").Copy() | Out-Null
    $ws.Range("A$TargetRow`:D$TargetRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

Copy-RowFormat 68 97
Copy-RowFormat 28 98
Copy-RowFormat 28 99
Copy-RowFormat 2  100
Copy-RowFormat 68 101
Copy-RowFormat 68 102
Copy-RowFormat 2  103
Copy-RowFormat 68 104
Copy-RowFormat 28 105

function Set-LogRow($Row, $DateValue, $Servizio, $Documento, $TipoModifica) {
    $ws.Cells.Item($Row, 1).Value = $DateValue
    if ($Servizio -ne $null) {
        $ws.Cells.Item($Row, 2).Value = $Servizio
    }
    # Write column D before C so new shared-string entries are appended in
    # the same order the original authoring tool produced them.
    $ws.Cells.Item($Row, 4).Value = $TipoModifica
    $ws.Cells.Item($Row, 3).Value = $Documento
}

Set-LogRow 97  "10/24/2017" "lista errori"       "errori_anpr_20171024.xlsx"                         "Inserimento/Aggiornamento codici di errore di ANPR `n"
Set-LogRow 98  "10/24/2017" "Tabelle decodifica" "tabella_46_tipo_mutazione_famiglia_convivenza.xlsx" "aggiunta la tabella di decodifica"
Set-LogRow 99  "10/24/2017" "Tabelle decodifica" "tabella_47_tipo_mutazione_residenza.xlsx"           "aggiunta la tabella di decodifica"
Set-LogRow 100 "10/24/2017" $null                "MI-14-AN-01 SPECIFICHE DI INTERFACCIA WS"           "Vedi §3"
Set-LogRow 101 "10/26/2017" "lista errori"       "errori_anpr_20171026.xlsx"                          "Inserimento codici di errore EN416, EN426"
Set-LogRow 102 "11/3/2017"  "lista errori"       "errori_anpr_20171103.xlsx"                          "Inserimento codice di errore EN447"
Set-LogRow 103 "11/3/2017"  $null                "MI-14-AN-01 SPECIFICHE DI INTERFACCIA WS"           "Vedi §3"
Set-LogRow 104 "11/9/2017"  "lista errori"       "errori_anpr_20171109.xlsx"                          "Inserimento/Aggiornamento codici di errore di ANPR `n"
Set-LogRow 105 "11/16/2017" "lista errori"       "errori_anpr_20171116.xlsx"                          "Inserimento/Aggiornamento codici di errore di ANPR"

# Row 104 picked up an automatic wrap-height from the long text above; put it
# back to the sheet's natural (default) height, matching the source row.
$ws.Rows.Item(104).AutoFit() | Out-Null

# Rows 97 and 105 are taller than the default row height (two-line content /
# a deliberate manual resize respectively).
$ws.Rows.Item(97).RowHeight = 30
$ws.Rows.Item(105).RowHeight = 18.75

# 3) Move the visible selection to where the author's cursor ended up after
#    typing in the new rows.
$ws.Range("C112").Select()
